$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.147438287734985
$ws.Range("B1").Value = 2.954680442810059
$ws.Range("C1").Value = 3.695252180099487
$ws.Range("D1").Value = 3.586387872695923
$ws.Range("E1").Value = 1.200913071632385
